$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.002767164346438711
$ws.Range("D2").Value = 0.005534328692877489
$ws.Range("E2").Value = 0.002767164346438755
$ws.Range("F2").Value = 0.008301493039316252
$ws.Range("G2").Value = 0.005534328692877542
$ws.Range("H2").Value = 0.008301493039316252
$ws.Range("I2").Value = 0.0166029860786326
$ws.Range("J2").Value = 0.0110686573857551

$ws.Range("C3").Value = 0.001396708916112762
$ws.Range("D3").Value = 0.0027934178322257
$ws.Range("E3").Value = 0.001396708916112888
$ws.Range("F3").Value = 0.004190126748338587
$ws.Range("G3").Value = 0.002793417832225825
$ws.Range("H3").Value = 0.004190126748338714
$ws.Range("I3").Value = 0.008380253496677301
$ws.Range("J3").Value = 0.005586835664451475

$ws.Range("C4").Value = 0.01709710209335594
$ws.Range("D4").Value = 0.03419420418671198
$ws.Range("E4").Value = 0.01709710209335608
$ws.Range("F4").Value = 0.05129130628006797
$ws.Range("G4").Value = 0.0341942041867121
$ws.Range("H4").Value = 0.05129130628006809
$ws.Range("I4").Value = 0.1025826125601361
$ws.Range("J4").Value = 0.06838840837342415

$ws.Range("C5").Value = 0.008776272009646518
$ws.Range("D5").Value = 0.0175525440192929
$ws.Range("E5").Value = 0.008776272009646428
$ws.Range("F5").Value = 0.02632881602893943
$ws.Range("G5").Value = 0.01755254401929289
$ws.Range("H5").Value = 0.02632881602893932
$ws.Range("I5").Value = 0.05265763205787874
$ws.Range("J5").Value = 0.03510508803858584

$ws.Range("C6").Value = 0.0045608954098372
$ws.Range("D6").Value = 0.009121790819674305
$ws.Range("E6").Value = 0.004560895409837224
$ws.Range("F6").Value = 0.0136826862295115
$ws.Range("G6").Value = 0.009121790819674274
$ws.Range("H6").Value = 0.01368268622951148
$ws.Range("I6").Value = 0.027365372459023
$ws.Range("J6").Value = 0.0182435816393487

$ws.Range("C7").Value = 0.004265792029994007
$ws.Range("D7").Value = 0.008531584059987994
$ws.Range("E7").Value = 0.004265792029993986
$ws.Range("F7").Value = 0.01279737608998207
$ws.Range("G7").Value = 0.008531584059988061
$ws.Range("H7").Value = 0.01279737608998205
$ws.Range("I7").Value = 0.02559475217996412
$ws.Range("J7").Value = 0.01706316811997605

$ws.Range("C8").Value = 0.02180821962234357
$ws.Range("D8").Value = 0.04361643924468724
$ws.Range("E8").Value = 0.02180821962234376
$ws.Range("F8").Value = 0.06542465886703092
$ws.Range("G8").Value = 0.04361643924468742
$ws.Range("H8").Value = 0.06542465886703108
$ws.Range("I8").Value = 0.130849317734062
$ws.Range("J8").Value = 0.08723287848937468

$ws.Range("C9").Value = 0.007324199000941048
$ws.Range("D9").Value = 0.01464839800188193
$ws.Range("E9").Value = 0.007324199000940878
$ws.Range("F9").Value = 0.02197259700282288
$ws.Range("G9").Value = 0.01464839800188183
$ws.Range("H9").Value = 0.02197259700282271
$ws.Range("I9").Value = 0.04394519400564559
$ws.Range("J9").Value = 0.02929679600376376

$ws.Range("C10").Value = 0.01388057970233389
$ws.Range("D10").Value = 0.02776115940466778
$ws.Range("E10").Value = 0.01388057970233388
$ws.Range("F10").Value = 0.04164173910700156
$ws.Range("G10").Value = 0.02776115940466768
$ws.Range("H10").Value = 0.04164173910700146
$ws.Range("I10").Value = 0.08328347821400302
$ws.Range("J10").Value = 0.05552231880933533

$ws.Range("C11").Value = 0.05155846797659717
$ws.Range("D11").Value = 0.1031169359531942
$ws.Range("E11").Value = 0.05155846797659698
$ws.Range("F11").Value = 0.1546754039297912
$ws.Range("G11").Value = 0.1031169359531941
$ws.Range("H11").Value = 0.1546754039297911
$ws.Range("I11").Value = 0.3093508078595822
$ws.Range("J11").Value = 0.2062338719063881

$ws.Range("C12").Value = 0.01162168521397276
$ws.Range("D12").Value = 0.02324337042794564
$ws.Range("E12").Value = 0.01162168521397285
$ws.Range("F12").Value = 0.03486505564191849
$ws.Range("G12").Value = 0.02324337042794579
$ws.Range("H12").Value = 0.03486505564191864
$ws.Range("I12").Value = 0.06973011128383713
$ws.Range("J12").Value = 0.0464867408558915

$ws.Range("C13").Value = 0.02431058797811773
$ws.Range("D13").Value = 0.0486211759562355
$ws.Range("E13").Value = 0.02431058797811765
$ws.Range("F13").Value = 0.07293176393435323
$ws.Range("G13").Value = 0.04862117595623548
$ws.Range("H13").Value = 0.07293176393435313
$ws.Range("I13").Value = 0.1458635278687063
$ws.Range("J13").Value = 0.09724235191247088

$ws.Range("C14").Value = 0.002435470696372353
$ws.Range("D14").Value = 0.004870941392744757
$ws.Range("E14").Value = 0.002435470696372404
$ws.Range("F14").Value = 0.007306412089117144
$ws.Range("G14").Value = 0.004870941392744792
$ws.Range("H14").Value = 0.007306412089117261
$ws.Range("I14").Value = 0.0146128241782344
$ws.Range("J14").Value = 0.009741882785489648

$ws.Range("C15").Value = 0.003769656843364089
$ws.Range("D15").Value = 0.007539313686728323
$ws.Range("E15").Value = 0.003769656843364123
$ws.Range("F15").Value = 0.01130897053009239
$ws.Range("G15").Value = 0.007539313686728309
$ws.Range("H15").Value = 0.01130897053009243
$ws.Range("I15").Value = 0.02261794106018482
$ws.Range("J15").Value = 0.01507862737345652

$ws.Range("C16").Value = 0.02030988836669417
$ws.Range("D16").Value = 0.04061977673338836
$ws.Range("E16").Value = 0.0203098883666941
$ws.Range("F16").Value = 0.06092966510008245
$ws.Range("G16").Value = 0.04061977673338828
$ws.Range("H16").Value = 0.06092966510008238
$ws.Range("I16").Value = 0.1218593302001648
$ws.Range("J16").Value = 0.08123955346677655

$ws.Range("C17").Value = 0.0006876968899117854
$ws.Range("D17").Value = 0.001375393779823474
$ws.Range("E17").Value = 0.0006876968899117142
$ws.Range("F17").Value = 0.002063090669735162
$ws.Range("G17").Value = 0.001375393779823377
$ws.Range("H17").Value = 0.002063090669735091
$ws.Range("I17").Value = 0.00412618133947028
$ws.Range("J17").Value = 0.002750787559646904

$ws.Range("C18").Value = 0.007303228522353748
$ws.Range("D18").Value = 0.01460645704470749
$ws.Range("E18").Value = 0.007303228522353616
$ws.Range("F18").Value = 0.02190968556706121
$ws.Range("G18").Value = 0.01460645704470741
$ws.Range("H18").Value = 0.02190968556706121
$ws.Range("I18").Value = 0.04381937113412236
$ws.Range("J18").Value = 0.02921291408941489

$ws.Range("C19").Value = 0.001905751655668409
$ws.Range("D19").Value = 0.003811503311336879
$ws.Range("E19").Value = 0.00190575165566834
$ws.Range("F19").Value = 0.005717254967005185
$ws.Range("G19").Value = 0.003811503311336732
$ws.Range("H19").Value = 0.005717254967005099
$ws.Range("I19").Value = 0.01143450993401026
$ws.Range("J19").Value = 0.007623006622673508

$ws.Range("C20").Value = 0.01241575855806621
$ws.Range("D20").Value = 0.02483151711613242
$ws.Range("E20").Value = 0.01241575855806614
$ws.Range("F20").Value = 0.03724727567419856
$ws.Range("G20").Value = 0.02483151711613235
$ws.Range("H20").Value = 0.03724727567419848
$ws.Range("I20").Value = 0.07449455134839704
$ws.Range("J20").Value = 0.0496630342322647

$ws.Range("C21").Value = 0.01924993406000224
$ws.Range("D21").Value = 0.03849986812000462
$ws.Range("E21").Value = 0.0192499340600023
$ws.Range("F21").Value = 0.05774980218000685
$ws.Range("G21").Value = 0.03849986812000461
$ws.Range("H21").Value = 0.05774980218000691
$ws.Range("I21").Value = 0.1154996043600138
$ws.Range("J21").Value = 0.07699973624000915

$ws.Range("C22").Value = 0.004960025362020695
$ws.Range("D22").Value = 0.009920050724041308
$ws.Range("E22").Value = 0.004960025362020604
$ws.Range("F22").Value = 0.01488007608606189
$ws.Range("G22").Value = 0.009920050724041217
$ws.Range("H22").Value = 0.0148800760860618
$ws.Range("I22").Value = 0.02976015217212371
$ws.Range("J22").Value = 0.01984010144808251

$ws.Range("C23").Value = 0.005069500098151526
$ws.Range("D23").Value = 0.01013900019630322
$ws.Range("E23").Value = 0.005069500098151639
$ws.Range("F23").Value = 0.01520850029445485
$ws.Range("G23").Value = 0.0101390001963033
$ws.Range("H23").Value = 0.01520850029445494
$ws.Range("I23").Value = 0.03041700058890979
$ws.Range("J23").Value = 0.02027800039260657

$ws.Range("C24").Value = 0.01483409318935364
$ws.Range("D24").Value = 0.02966818637870732
$ws.Range("E24").Value = 0.01483409318935368
$ws.Range("F24").Value = 0.04450227956806091
$ws.Range("G24").Value = 0.02966818637870721
$ws.Range("H24").Value = 0.04450227956806089
$ws.Range("I24").Value = 0.08900455913612178
$ws.Range("J24").Value = 0.05933637275741453

$ws.Range("C25").Value = 0.00692772811322721
$ws.Range("D25").Value = 0.0138554562264545
$ws.Range("E25").Value = 0.006927728113227288
$ws.Range("F25").Value = 0.02078318433968171
$ws.Range("G25").Value = 0.0138554562264545
$ws.Range("H25").Value = 0.02078318433968179
$ws.Range("I25").Value = 0.0415663686793635
$ws.Range("J25").Value = 0.027710912452909

